$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 57.8759039791103
$ws.Range("B2").Value = -27907.9253661154
$ws.Range("C2").Value = 55934.2673103241
$ws.Range("D2").Value = 56391.6984294377
$ws.Range("E2").Value = 1244.83295386678
$ws.Range("F2").Value = 16686.1240960209
